$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 794.1667
$ws.Range("I12").Value = 753.3
$ws.Range("K12").Value = 753.3
$ws.Range("M12").Value = -583.3

# Row 51
$ws.Range("H51").Value = 3204.7058
$ws.Range("I51").Value = 2993.5
$ws.Range("J51").Value = 3232.8667
$ws.Range("K51").Value = 2993.5
$ws.Range("L51").Value = 3232.8667
$ws.Range("M51").Value = -2509.5
$ws.Range("N51").Value = -4200.8667

# Row 74
$ws.Range("H74").Value = 9700
$ws.Range("I74").Value = 9700
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 9700
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -8764
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 9700
$ws.Range("I77").Value = 9700
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 48500
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -43820
$ws.Range("N77").ClearContents()

# Row 112
$ws.Range("H112").Value = 2171.6667
$ws.Range("J112").Value = 1382.5
$ws.Range("L112").Value = 4147.5
$ws.Range("N112").Value = -6363.5

# Row 132
$ws.Range("H132").Value = 3455.5
$ws.Range("I132").Value = 1213.2667
$ws.Range("K132").Value = 3639.800099999999
$ws.Range("M132").Value = -1109.800099999999

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3116.2354
$ws.Range("I45").Value = 1975.6666
$ws.Range("J45").Value = 3738.3635
$ws.Range("K45").Value = 1975.6666
$ws.Range("L45").Value = 3738.3635
$ws.Range("M45").Value = -1598.6666
$ws.Range("N45").Value = -4492.363499999999

# Row 61
$ws.Range("H61").Value = 2050.7334
$ws.Range("I61").Value = 1563.5
$ws.Range("K61").Value = 1563.5
$ws.Range("M61").Value = -1351.5

# Row 74
$ws.Range("H74").Value = 3417.2727
$ws.Range("I74").Value = 2662.889
$ws.Range("K74").Value = 2662.889
$ws.Range("M74").Value = -1788.889

# Row 77
$ws.Range("H77").Value = 3417.2727
$ws.Range("I77").Value = 2662.889
$ws.Range("K77").Value = 13314.445
$ws.Range("M77").Value = -8946.445

# Row 110
$ws.Range("H110").Value = 2441.8
$ws.Range("I110").Value = 2299
$ws.Range("K110").Value = 2299
$ws.Range("M110").Value = -254

# Row 122
$ws.Range("H122").Value = 3010
$ws.Range("I122").Value = 3010
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9030
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6580
$ws.Range("N122").ClearContents()

# Row 135
$ws.Range("H135").Value = 79997
$ws.Range("J135").Value = 79997
$ws.Range("L135").Value = 79997
$ws.Range("N135").Value = -90137

# Row 136
$ws.Range("H136").Value = 2050.7334
$ws.Range("I136").Value = 1563.5
$ws.Range("K136").Value = 4690.5
$ws.Range("M136").Value = -2140.5

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 27511.223
$ws.Range("I82").Value = 7857.8335
$ws.Range("J82").Value = 66818
$ws.Range("K82").Value = 7857.8335
$ws.Range("L82").Value = 66818
$ws.Range("M82").Value = -7474.8335
$ws.Range("N82").Value = -67584

# Row 85
$ws.Range("H85").Value = 27511.223
$ws.Range("I85").Value = 7857.8335
$ws.Range("J85").Value = 66818
$ws.Range("K85").Value = 7857.8335
$ws.Range("L85").Value = 66818
$ws.Range("M85").Value = -6531.8335
$ws.Range("N85").Value = -69470

# Row 86
$ws.Range("H86").Value = 4298.9
$ws.Range("J86").Value = 5427
$ws.Range("L86").Value = 5427
$ws.Range("N86").Value = -7673

# Row 89
$ws.Range("H89").Value = 4298.9
$ws.Range("J89").Value = 5427
$ws.Range("L89").Value = 27135
$ws.Range("N89").Value = -38367

# Row 99
$ws.Range("H99").Value = 1199.6
$ws.Range("I99").Value = 1376
$ws.Range("J99").Value = 494
$ws.Range("K99").Value = 1376
$ws.Range("L99").Value = 494
$ws.Range("M99").Value = 122
$ws.Range("N99").Value = -3490

# Row 107
$ws.Range("H107").Value = 9175.25
$ws.Range("I107").Value = 9201.714
$ws.Range("K107").Value = 9201.714
$ws.Range("M107").Value = -7281.714

# Row 135
$ws.Range("H135").Value = 186666.33
$ws.Range("J135").Value = 186666.33
$ws.Range("L135").Value = 186666.33
$ws.Range("N135").Value = -196806.33

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 3743.7334
$ws.Range("I134").Value = 2420.7273
$ws.Range("K134").Value = 7262.1819
$ws.Range("M134").Value = -4727.1819

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5999.5
$ws.Range("I70").Value = 5999.5
$ws.Range("K70").Value = 5999.5
$ws.Range("M70").Value = -5729.5

# Row 73
$ws.Range("H73").Value = 5999.5
$ws.Range("I73").Value = 5999.5
$ws.Range("K73").Value = 5999.5
$ws.Range("M73").Value = -5063.5

# Row 80
$ws.Range("H80").Value = 1159.4
$ws.Range("I80").Value = 1199
$ws.Range("J80").Value = 1149.5
$ws.Range("K80").Value = 1199
$ws.Range("L80").Value = 1149.5
$ws.Range("M80").Value = -201
$ws.Range("N80").Value = -3145.5

# Row 83
$ws.Range("H83").Value = 1159.4
$ws.Range("I83").Value = 1199
$ws.Range("J83").Value = 1149.5
$ws.Range("K83").Value = 5995
$ws.Range("L83").Value = 5747.5
$ws.Range("M83").Value = -1003
$ws.Range("N83").Value = -15731.5

# Row 102
$ws.Range("H102").Value = 3508.5
$ws.Range("I102").Value = 3253.0833
$ws.Range("J102").Value = 4274.75
$ws.Range("K102").Value = 3253.0833
$ws.Range("L102").Value = 4274.75
$ws.Range("M102").Value = -1631.0833
$ws.Range("N102").Value = -7518.75

# Row 113
$ws.Range("H113").Value = 6381
$ws.Range("I113").Value = 2787
$ws.Range("K113").Value = 2787
$ws.Range("M113").Value = -617

# Row 122
$ws.Range("H122").Value = 3738.0715
$ws.Range("I122").Value = 3645
$ws.Range("K122").Value = 10935
$ws.Range("M122").Value = -8485

# Row 126
$ws.Range("H126").Value = 2500
$ws.Range("J126").Value = 2500
$ws.Range("L126").Value = 7500
$ws.Range("N126").Value = -12440

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2037.5
$ws.Range("I16").Value = 1445.2
$ws.Range("J16").Value = 4999
$ws.Range("K16").Value = 1445.2
$ws.Range("L16").Value = 4999
$ws.Range("M16").Value = -1275.2
$ws.Range("N16").Value = -5339

# Row 22
$ws.Range("H22").Value = 1548.25
$ws.Range("I22").Value = 1177.2
$ws.Range("J22").Value = 2166.6667
$ws.Range("K22").Value = 1177.2
$ws.Range("L22").Value = 2166.6667
$ws.Range("M22").Value = -882.2
$ws.Range("N22").Value = -2756.6667

# Row 27
$ws.Range("H27").Value = 1548.25
$ws.Range("I27").Value = 1177.2
$ws.Range("J27").Value = 2166.6667
$ws.Range("K27").Value = 1177.2
$ws.Range("L27").Value = 2166.6667
$ws.Range("M27").Value = -1070.2
$ws.Range("N27").Value = -2380.6667

# Row 40
$ws.Range("H40").Value = 6135.148
$ws.Range("I40").Value = 5397.696
$ws.Range("J40").Value = 10375.5
$ws.Range("K40").Value = 5397.696
$ws.Range("L40").Value = 10375.5
$ws.Range("M40").Value = -5261.696
$ws.Range("N40").Value = -10647.5

# Row 82
$ws.Range("H82").Value = 3874.4546
$ws.Range("I82").Value = 650
$ws.Range("J82").Value = 4591
$ws.Range("K82").Value = 650
$ws.Range("L82").Value = 4591
$ws.Range("M82").Value = -289
$ws.Range("N82").Value = -5313

# Row 85
$ws.Range("H85").Value = 3874.4546
$ws.Range("I85").Value = 650
$ws.Range("J85").Value = 4591
$ws.Range("K85").Value = 650
$ws.Range("L85").Value = 4591
$ws.Range("M85").Value = 598
$ws.Range("N85").Value = -7087

# Row 101
$ws.Range("H101").Value = 11472.4
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 11472.4
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 11472.4
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -17962.4

# Row 122
$ws.Range("H122").Value = 1704
$ws.Range("I122").Value = 1704
$ws.Range("K122").Value = 5112
$ws.Range("M122").Value = -2662

$ws = $wb.Worksheets.Item("WVR")
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

# Row 100
$ws.Range("H100").Value = 1399
$ws.Range("I100").Value = 1399
$ws.Range("K100").Value = 2798
$ws.Range("M100").Value = -2257

# Row 122
$ws.Range("H122").Value = 1779
$ws.Range("I122").Value = 1549.3334
$ws.Range("K122").Value = 4648.0002
$ws.Range("M122").Value = -2198.0002
